$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 341.54544
$ws.Range("I53").Value = 312.83334
$ws.Range("K53").Value = 312.83334
$ws.Range("M53").Value = 324.16666
$ws.Range("H93").Value = 44937.5
$ws.Range("J93").Value = 44937.5
$ws.Range("L93").Value = 44937.5
$ws.Range("N93").Value = -49929.5
$ws.Range("H107").Value = 338.83334
$ws.Range("I107").Value = 184.375
$ws.Range("J107").Value = 462.4
$ws.Range("K107").Value = 184.375
$ws.Range("L107").Value = 462.4
$ws.Range("M107").Value = 1735.625
$ws.Range("N107").Value = -4302.4
$ws.Range("H116").Value = 2553.85
$ws.Range("I116").Value = 1657
$ws.Range("J116").Value = 3036.7693
$ws.Range("K116").Value = 1657
$ws.Range("L116").Value = 3036.7693
$ws.Range("M116").Value = 1785
$ws.Range("N116").Value = -9920.7693
$ws.Range("H131").Value = 2298
$ws.Range("I131").Value = 999
$ws.Range("J131").Value = 2947.5
$ws.Range("K131").Value = 2997
$ws.Range("L131").Value = 8842.5
$ws.Range("M131").Value = 2043
$ws.Range("N131").Value = -18922.5
$ws.Range("H132").Value = 29616.648
$ws.Range("I132").Value = 34923.16
$ws.Range("J132").Value = 2199.6667
$ws.Range("K132").Value = 104769.48
$ws.Range("L132").Value = 6599.000100000001
$ws.Range("M132").Value = -102239.48
$ws.Range("N132").Value = -11659.0001
$ws.Range("H138").Value = 2958.8684
$ws.Range("I138").Value = 2256.0833
$ws.Range("J138").Value = 3283.2307
$ws.Range("K138").Value = 6768.249899999999
$ws.Range("L138").Value = 9849.6921
$ws.Range("M138").Value = -1628.249899999999
$ws.Range("N138").Value = -20129.6921
$ws.Range("H140").Value = 49800
$ws.Range("J140").Value = 49800
$ws.Range("L140").Value = 49800
$ws.Range("N140").Value = -60160

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 634.1667
$ws.Range("I5").Value = 738.75
$ws.Range("J5").Value = 425
$ws.Range("K5").Value = 738.75
$ws.Range("L5").Value = 425
$ws.Range("M5").Value = -626.75
$ws.Range("N5").Value = -649
$ws.Range("H32").Value = 2127.8657
$ws.Range("I32").Value = 1409
$ws.Range("J32").Value = 7429.5
$ws.Range("K32").Value = 1409
$ws.Range("L32").Value = 7429.5
$ws.Range("M32").Value = -1122
$ws.Range("N32").Value = -8003.5
$ws.Range("H38").Value = 7669.3335
$ws.Range("I38").Value = 7669.3335
$ws.Range("K38").Value = 7669.3335
$ws.Range("M38").Value = -7202.3335
$ws.Range("H63").Value = 2234314.5
$ws.Range("I63").Value = 2117
$ws.Range("J63").Value = 15627500
$ws.Range("K63").Value = 2117
$ws.Range("L63").Value = 15627500
$ws.Range("M63").Value = -1431
$ws.Range("N63").Value = -15628872
$ws.Range("H66").Value = 2234314.5
$ws.Range("I66").Value = 2117
$ws.Range("J66").Value = 15627500
$ws.Range("K66").Value = 10585
$ws.Range("L66").Value = 78137500
$ws.Range("M66").Value = -7153
$ws.Range("N66").Value = -78144364
$ws.Range("H74").Value = 651.0833
$ws.Range("I74").Value = 310.79166
$ws.Range("J74").Value = 1331.6666
$ws.Range("K74").Value = 310.79166
$ws.Range("L74").Value = 1331.6666
$ws.Range("M74").Value = 563.20834
$ws.Range("N74").Value = -3079.6666
$ws.Range("H77").Value = 651.0833
$ws.Range("I77").Value = 310.79166
$ws.Range("J77").Value = 1331.6666
$ws.Range("K77").Value = 1553.9583
$ws.Range("L77").Value = 6658.333000000001
$ws.Range("M77").Value = 2814.0417
$ws.Range("N77").Value = -15394.333
$ws.Range("H110").Value = 2448.3333
$ws.Range("I110").Value = 2448.3333
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 2448.3333
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = -403.3332999999998
$ws.Range("N110").ClearContents()
$ws.Range("H122").Value = 1580.1666
$ws.Range("I122").Value = 1518.4348
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 4555.3044
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -2105.3044
$ws.Range("N122").Value = -13900

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 634.1667
$ws.Range("I4").Value = 738.75
$ws.Range("J4").Value = 425
$ws.Range("K4").Value = 738.75
$ws.Range("L4").Value = 425
$ws.Range("M4").Value = -623.75
$ws.Range("N4").Value = -655
$ws.Range("H105").Value = 2382630.2
$ws.Range("I105").Value = 1723.4
$ws.Range("J105").Value = 4547091
$ws.Range("K105").Value = 1723.4
$ws.Range("L105").Value = 4547091
$ws.Range("M105").Value = 23.59999999999991
$ws.Range("N105").Value = -4550585
$ws.Range("H107").Value = 1443.5
$ws.Range("I107").Value = 1332.2
$ws.Range("J107").Value = 2000
$ws.Range("K107").Value = 1332.2
$ws.Range("L107").Value = 2000
$ws.Range("M107").Value = 587.8
$ws.Range("N107").Value = -5840

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 674.8333
$ws.Range("I22").Value = 712.25
$ws.Range("J22").Value = 600
$ws.Range("K22").Value = 712.25
$ws.Range("L22").Value = 600
$ws.Range("M22").Value = -362.25
$ws.Range("N22").Value = -1300
$ws.Range("H94").Value = 5528.4443
$ws.Range("J94").Value = 7651.2
$ws.Range("L94").Value = 7651.2
$ws.Range("N94").Value = -8553.200000000001
$ws.Range("H107").Value = 1772.3334
$ws.Range("I107").Value = 2126.3333
$ws.Range("J107").Value = 1418.3334
$ws.Range("K107").Value = 2126.3333
$ws.Range("L107").Value = 1418.3334
$ws.Range("M107").Value = -206.3332999999998
$ws.Range("N107").Value = -5258.3334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 927.2857
$ws.Range("I68").Value = 622.2222
$ws.Range("J68").Value = 1250.2941
$ws.Range("K68").Value = 1866.6666
$ws.Range("L68").Value = 3750.8823
$ws.Range("M68").Value = -1055.6666
$ws.Range("N68").Value = -5372.8823
$ws.Range("H71").Value = 927.2857
$ws.Range("I71").Value = 622.2222
$ws.Range("J71").Value = 1250.2941
$ws.Range("K71").Value = 5599.999800000001
$ws.Range("L71").Value = 11252.6469
$ws.Range("M71").Value = -1543.999800000001
$ws.Range("N71").Value = -19364.6469
$ws.Range("H113").Value = 311.55554
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 311.55554
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 934.66662
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -5274.66662
$ws.Range("H131").Value = 802.8200000000001
$ws.Range("J131").Value = 823.0625
$ws.Range("L131").Value = 2469.1875
$ws.Range("N131").Value = -12549.1875
$ws.Range("H132").Value = 1911
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 1911
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 17199
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -22259
$ws.Range("H137").Value = 2920.5625
$ws.Range("J137").Value = 3266.3572
$ws.Range("L137").Value = 9799.071599999999
$ws.Range("N137").Value = -19999.0716

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3188.7827
$ws.Range("I80").Value = 2850
$ws.Range("J80").Value = 3449.3845
$ws.Range("K80").Value = 2850
$ws.Range("L80").Value = 3449.3845
$ws.Range("M80").Value = -1852
$ws.Range("N80").Value = -5445.3845
$ws.Range("H83").Value = 3188.7827
$ws.Range("I83").Value = 2850
$ws.Range("J83").Value = 3449.3845
$ws.Range("K83").Value = 14250
$ws.Range("L83").Value = 17246.9225
$ws.Range("M83").Value = -9258
$ws.Range("N83").Value = -27230.9225
$ws.Range("H107").Value = 1341.091
$ws.Range("I107").Value = 369.125
$ws.Range("J107").Value = 3933
$ws.Range("K107").Value = 369.125
$ws.Range("L107").Value = 3933
$ws.Range("M107").Value = 1550.875
$ws.Range("N107").Value = -7773
$ws.Range("H113").Value = 2664.8823
$ws.Range("I113").Value = 2615.3845
$ws.Range("J113").Value = 2825.75
$ws.Range("K113").Value = 2615.3845
$ws.Range("L113").Value = 2825.75
$ws.Range("M113").Value = -445.3845000000001
$ws.Range("N113").Value = -7165.75
$ws.Range("H132").Value = 28231.762
$ws.Range("I132").Value = 5613.3335
$ws.Range("J132").Value = 58389.668
$ws.Range("K132").Value = 16840.0005
$ws.Range("L132").Value = 175169.004
$ws.Range("M132").Value = -14310.0005
$ws.Range("N132").Value = -180229.004

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1155
$ws.Range("I16").Value = 1186
$ws.Range("J16").Value = 1000
$ws.Range("K16").Value = 1186
$ws.Range("L16").Value = 1000
$ws.Range("M16").Value = -1016
$ws.Range("N16").Value = -1340
$ws.Range("H68").Value = 4888.2
$ws.Range("I68").Value = 3496.1
$ws.Range("J68").Value = 6280.3
$ws.Range("K68").Value = 3496.1
$ws.Range("L68").Value = 6280.3
$ws.Range("M68").Value = -2747.1
$ws.Range("N68").Value = -7778.3
$ws.Range("H71").Value = 4888.2
$ws.Range("I71").Value = 3496.1
$ws.Range("J71").Value = 6280.3
$ws.Range("K71").Value = 17480.5
$ws.Range("L71").Value = 31401.5
$ws.Range("M71").Value = -13736.5
$ws.Range("N71").Value = -38889.5
$ws.Range("H82").Value = 1971.0834
$ws.Range("I82").Value = 2664.7144
$ws.Range("K82").Value = 2664.7144
$ws.Range("M82").Value = -2303.7144
$ws.Range("H85").Value = 1971.0834
$ws.Range("I85").Value = 2664.7144
$ws.Range("K85").Value = 2664.7144
$ws.Range("M85").Value = -1416.7144

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5142.857
$ws.Range("J62").Value = 5142.857
$ws.Range("L62").Value = 5142.857
$ws.Range("N62").Value = -6390.857
$ws.Range("H65").Value = 5142.857
$ws.Range("J65").Value = 5142.857
$ws.Range("L65").Value = 25714.285
$ws.Range("N65").Value = -31954.285
$ws.Range("H107").Value = 1983.5555
$ws.Range("I107").Value = 2101
$ws.Range("J107").Value = 1950
$ws.Range("K107").Value = 6303
$ws.Range("L107").Value = 5850
$ws.Range("M107").Value = -4383
$ws.Range("N107").Value = -9690
$ws.Range("H126").Value = 1082.7142
$ws.Range("I126").Value = 955.8
$ws.Range("J126").Value = 1400
$ws.Range("K126").Value = 2867.4
$ws.Range("L126").Value = 4200
$ws.Range("M126").Value = -397.3999999999996
$ws.Range("N126").Value = -9140
$ws.Range("H132").Value = 1883.2858
$ws.Range("I132").Value = 1385.3529
$ws.Range("J132").Value = 3999.5
$ws.Range("K132").Value = 4156.0587
$ws.Range("L132").Value = 11998.5
$ws.Range("M132").Value = -1626.0587
$ws.Range("N132").Value = -17058.5
